$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("n1_d40")
$ws.Cells.Item(2, 6).Value = 0.0284566
$ws.Cells.Item(2, 7).Value = 5.6
$ws.Cells.Item(3, 6).Value = 0.0284523
$ws.Cells.Item(3, 7).Value = 5.6
$ws.Cells.Item(4, 6).Value = 0.02848
$ws.Cells.Item(4, 7).Value = 5.7
$ws.Cells.Item(5, 6).Value = 0.0284566
$ws.Cells.Item(5, 7).Value = 5.6
$ws.Cells.Item(6, 6).Value = 0.0284692
$ws.Cells.Item(6, 7).Value = 5.7
$ws.Cells.Item(7, 6).Value = 0.0284842
$ws.Cells.Item(7, 7).Value = 5.6
$ws.Cells.Item(8, 6).Value = 0.0284033
$ws.Cells.Item(8, 7).Value = 5.6
$ws.Cells.Item(9, 6).Value = 0.0286137
$ws.Cells.Item(9, 7).Value = 5.6
$ws.Cells.Item(10, 6).Value = 0.028453
$ws.Cells.Item(10, 7).Value = 5.6
$ws.Cells.Item(11, 6).Value = 0.0284505
$ws.Cells.Item(11, 7).Value = 5.7
$ws.Cells.Item(12, 6).Value = 0.02847193999999999
$ws.Cells.Item(12, 7).Value = 5.630000000000001

$ws = $wb.Worksheets.Item("n1_d60")
$ws.Cells.Item(2, 6).Value = 0.0381096
$ws.Cells.Item(2, 7).Value = 8.6
$ws.Cells.Item(3, 6).Value = 0.0390949
$ws.Cells.Item(3, 7).Value = 8.699999999999999
$ws.Cells.Item(4, 6).Value = 0.039149
$ws.Cells.Item(4, 7).Value = 8.699999999999999
$ws.Cells.Item(5, 6).Value = 0.0381007
$ws.Cells.Item(5, 7).Value = 8.6
$ws.Cells.Item(6, 6).Value = 0.0381149
$ws.Cells.Item(6, 7).Value = 8.6
$ws.Cells.Item(7, 6).Value = 0.0402652
$ws.Cells.Item(7, 7).Value = 8.1
$ws.Cells.Item(8, 6).Value = 0.0379945
$ws.Cells.Item(8, 7).Value = 8.699999999999999
$ws.Cells.Item(9, 6).Value = 0.0393364
$ws.Cells.Item(9, 7).Value = 8.199999999999999
$ws.Cells.Item(10, 6).Value = 0.0380367
$ws.Cells.Item(10, 7).Value = 8.699999999999999
$ws.Cells.Item(11, 6).Value = 0.0391196
$ws.Cells.Item(11, 7).Value = 8.699999999999999
$ws.Cells.Item(12, 6).Value = 0.03873215
$ws.Cells.Item(12, 7).Value = 8.559999999999999

$ws = $wb.Worksheets.Item("n1_d80")
$ws.Cells.Item(2, 6).Value = 0.0476614
$ws.Cells.Item(2, 7).Value = 11.6
$ws.Cells.Item(3, 6).Value = 0.0477105
$ws.Cells.Item(3, 7).Value = 11.5
$ws.Cells.Item(4, 6).Value = 0.04763
$ws.Cells.Item(4, 7).Value = 11.5
$ws.Cells.Item(5, 6).Value = 0.0476167
$ws.Cells.Item(5, 7).Value = 11.4
$ws.Cells.Item(6, 6).Value = 0.0487558
$ws.Cells.Item(6, 7).Value = 11.5
$ws.Cells.Item(7, 6).Value = 0.0478049
$ws.Cells.Item(7, 7).Value = 11.6
$ws.Cells.Item(8, 6).Value = 0.0498273
$ws.Cells.Item(8, 7).Value = 11.5
$ws.Cells.Item(9, 6).Value = 0.0477081
$ws.Cells.Item(9, 7).Value = 11.6
$ws.Cells.Item(10, 6).Value = 0.0478189
$ws.Cells.Item(10, 7).Value = 11.4
$ws.Cells.Item(11, 6).Value = 0.0476825
$ws.Cells.Item(11, 7).Value = 11.5
$ws.Cells.Item(12, 6).Value = 0.04802161
$ws.Cells.Item(12, 7).Value = 11.51

$ws = $wb.Worksheets.Item("n1_d100")
$ws.Cells.Item(2, 6).Value = 0.058346
$ws.Cells.Item(2, 7).Value = 13.3
$ws.Cells.Item(3, 6).Value = 0.0572606
$ws.Cells.Item(3, 7).Value = 13.3
$ws.Cells.Item(4, 6).Value = 0.057294
$ws.Cells.Item(4, 7).Value = 13.3
$ws.Cells.Item(5, 6).Value = 0.0616626
$ws.Cells.Item(5, 7).Value = 13.2
$ws.Cells.Item(6, 6).Value = 0.0583347
$ws.Cells.Item(6, 7).Value = 13.2
$ws.Cells.Item(7, 6).Value = 0.0574581
$ws.Cells.Item(7, 7).Value = 13.2
$ws.Cells.Item(8, 6).Value = 0.0572373
$ws.Cells.Item(8, 7).Value = 13.3
$ws.Cells.Item(9, 6).Value = 0.0604873
$ws.Cells.Item(9, 7).Value = 13.3
$ws.Cells.Item(10, 6).Value = 0.058366
$ws.Cells.Item(10, 7).Value = 13.2
$ws.Cells.Item(11, 6).Value = 0.0604756
$ws.Cells.Item(11, 7).Value = 13.3
$ws.Cells.Item(12, 6).Value = 0.05869222
$ws.Cells.Item(12, 7).Value = 13.26
